$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the FP_Card_No / FP_Pin values per the commit's edits.
# Column A holds long digit strings stored as text (shared strings) in the
# original file (a quote-prefixed numeric-format cell). Prefix the value
# with a leading apostrophe so Excel keeps storing it as text instead of
# auto-converting the digit string to a number.
$ws.Range("A2").Value = "'6375004101502496"
$ws.Range("B2").Value = 959

$ws.Range("A3").Value = "'6375004102003502"
$ws.Range("B3").Value = 571

$ws.Range("A4").Value = "'6375004101502496"
$ws.Range("B4").Value = 959

$ws.Range("A5").Value = "'6375004110391584"
$ws.Range("B5").Value = 779

$ws.Range("A6").Value = "'6375004102254139"
$ws.Range("B6").Value = 257

# Move the active selection to D6 (matches the saved selection in the file)
$ws.Range("D6").Select()
